$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "RM 232" (row 26).
$ws.Rows.Item(26).Delete()

# After the above deletion, the row that contained "SC 92" (originally row 28)
# has shifted up to row 27. Delete it too.
$ws.Rows.Item(27).Delete()

# Two cells in column F that were previously blank now have values.
$ws.Range("F27").Value = 17
$ws.Range("F30").Value = 16.89

# Three cells in column F that previously had values are now blank.
$ws.Range("F28").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("F32").ClearContents()
